$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.457.60"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "1.871.08"
$ws.Range("E3").Value = "  +1.60%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.54%  "

$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.87"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.09%  "

$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "2.135.48"

$ws.Range("D13").Value = "1.916.17"
$ws.Range("E13").Value = "  +4.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.63"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.685"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.10%  "

$ws.Range("D17").Value = "35.424.06"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.40"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "242.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +25.24%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0564"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.06"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.71%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +23.16%  "

$ws.Range("E35").Value = "  +8.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.823"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +18.30%  "

$ws.Range("E37").Value = "  +6.19%  "

$ws.Range("E38").Value = "  +4.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0204"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.43"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("D41").Value = "1.354.20"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.26"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.11%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0601"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +14.89%  "

$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.29"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +58.08%  "

$ws.Range("E45").Value = "  +2.90%  "

$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.66"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("D49").Value = "2.050.75"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("E50").Value = "  +3.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.39%  "
